$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2889.6562
$ws.Range("I64").Value = 3389.8333
$ws.Range("J64").Value = 2774.2307
$ws.Range("K64").Value = 3389.8333
$ws.Range("L64").Value = 2774.2307
$ws.Range("M64").Value = -3141.8333
$ws.Range("N64").Value = -3270.2307
$ws.Range("H67").Value = 2889.6562
$ws.Range("I67").Value = 3389.8333
$ws.Range("J67").Value = 2774.2307
$ws.Range("K67").Value = 3389.8333
$ws.Range("L67").Value = 2774.2307
$ws.Range("M67").Value = -2531.8333
$ws.Range("N67").Value = -4490.2307
$ws.Range("H113").Value = 2533.9033
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 2542.04
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 2542.04
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -9050.040000000001
$ws.Range("H116").Value = 2077.3125
$ws.Range("I116").Value = 1987.9269
$ws.Range("J116").Value = 2600.8572
$ws.Range("K116").Value = 1987.9269
$ws.Range("L116").Value = 2600.8572
$ws.Range("M116").Value = 1454.0731
$ws.Range("N116").Value = -9484.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2197.2856
$ws.Range("I45").Value = 2432.75
$ws.Range("J45").Value = 1883.3334
$ws.Range("K45").Value = 2432.75
$ws.Range("L45").Value = 1883.3334
$ws.Range("M45").Value = -2055.75
$ws.Range("N45").Value = -2637.3334
$ws.Range("H74").Value = 3657.3215
$ws.Range("I74").Value = 3187
$ws.Range("K74").Value = 3187
$ws.Range("M74").Value = -2313
$ws.Range("H77").Value = 3657.3215
$ws.Range("I77").Value = 3187
$ws.Range("K77").Value = 15935
$ws.Range("M77").Value = -11567
$ws.Range("H102").Value = 1140.6923
$ws.Range("I102").Value = 1152.5
$ws.Range("J102").Value = 999
$ws.Range("K102").Value = 1152.5
$ws.Range("L102").Value = 999
$ws.Range("M102").Value = 469.5
$ws.Range("N102").Value = -4243
$ws.Range("H110").Value = 2682.6667
$ws.Range("I110").Value = 2616.1428
$ws.Range("J110").Value = 2915.5
$ws.Range("K110").Value = 2616.1428
$ws.Range("L110").Value = 2915.5
$ws.Range("M110").Value = -571.1428000000001
$ws.Range("N110").Value = -7005.5
$ws.Range("H122").Value = 1870.4615
$ws.Range("I122").Value = 1130.5714
$ws.Range("J122").Value = 2733.6667
$ws.Range("K122").Value = 3391.7142
$ws.Range("L122").Value = 8201.000100000001
$ws.Range("M122").Value = -941.7142000000003
$ws.Range("N122").Value = -13101.0001
$ws.Range("H123").Value = 50799.8
$ws.Range("J123").Value = 50799.8
$ws.Range("L123").Value = 50799.8
$ws.Range("N123").Value = -60599.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7279.5806
$ws.Range("I94").Value = 780.625
$ws.Range("J94").Value = 29561.715
$ws.Range("K94").Value = 780.625
$ws.Range("L94").Value = 29561.715
$ws.Range("M94").Value = -329.625
$ws.Range("N94").Value = -30463.715
$ws.Range("H99").Value = 1503.5454
$ws.Range("I99").Value = 1421.5385
$ws.Range("J99").Value = 1622
$ws.Range("K99").Value = 1421.5385
$ws.Range("L99").Value = 1622
$ws.Range("M99").Value = 76.46149999999989
$ws.Range("N99").Value = -4618
$ws.Range("H105").Value = 5738.067
$ws.Range("I105").Value = 1800
$ws.Range("K105").Value = 1800
$ws.Range("M105").Value = -53

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1644.5
$ws.Range("I102").Value = 1692.6666
$ws.Range("K102").Value = 1692.6666
$ws.Range("M102").Value = -70.66660000000002
$ws.Range("H113").Value = 1006.5714
$ws.Range("I113").Value = 845.7778
$ws.Range("J113").Value = 1296
$ws.Range("K113").Value = 845.7778
$ws.Range("L113").Value = 1296
$ws.Range("M113").Value = 1324.2222
$ws.Range("N113").Value = -5636
$ws.Range("H122").Value = 1647.8182
$ws.Range("I122").Value = 1835.3334
$ws.Range("J122").Value = 804
$ws.Range("K122").Value = 5506.0002
$ws.Range("L122").Value = 2412
$ws.Range("M122").Value = -3056.0002
$ws.Range("N122").Value = -7312
$ws.Range("H123").Value = 12958.125
$ws.Range("J123").Value = 12958.125
$ws.Range("L123").Value = 12958.125
$ws.Range("N123").Value = -17858.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2980
$ws.Range("J3").Value = 2980
$ws.Range("L3").Value = 2980
$ws.Range("N3").Value = -3204
$ws.Range("H14").Value = 4500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 4500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 4500
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = -4844
$ws.Range("H15").Value = 2980
$ws.Range("J15").Value = 2980
$ws.Range("L15").Value = 2980
$ws.Range("N15").Value = -3320
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""
$ws.Range("H61").Value = 2689.375
$ws.Range("I61").Value = 2930.7144
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 2930.7144
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -2728.7144
$ws.Range("N61").Value = -1404
$ws.Range("H113").Value = 2689.375
$ws.Range("I113").Value = 2930.7144
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2930.7144
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -760.7143999999998
$ws.Range("N113").Value = -5340
$ws.Range("H122").Value = 2658.12
$ws.Range("I122").Value = 2660.8
$ws.Range("J122").Value = 2656.3333
$ws.Range("K122").Value = 7982.400000000001
$ws.Range("L122").Value = 7968.999899999999
$ws.Range("M122").Value = -5532.400000000001
$ws.Range("N122").Value = -12868.9999
$ws.Range("H132").Value = 9096729
$ws.Range("I132").Value = 20010172
$ws.Range("J132").Value = 2193.5667
$ws.Range("K132").Value = 60030516
$ws.Range("L132").Value = 6580.7001
$ws.Range("M132").Value = -60027986
$ws.Range("N132").Value = -11640.7001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1716.5
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1759.8
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 5279.4
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -10179.4
$ws.Range("H123").Value = 34115
$ws.Range("J123").Value = 34115
$ws.Range("L123").Value = 34115
$ws.Range("N123").Value = -43915
